# Insert a new weekly price-report row above row 248 (Apio, Chillán sheet).
# This shifts the existing rows 248-294 down to 249-295, preserving their
# data and formatting, and the new row 248 is populated with this week's
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(248).Insert()

$ws.Cells.Item(248, 1).Value  = 7
$ws.Cells.Item(248, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(248, 3).Value  = "Ñuble"
$ws.Cells.Item(248, 4).Value  = 44995
$ws.Cells.Item(248, 5).Value  = 16
$ws.Cells.Item(248, 6).Value  = 100112017
$ws.Cells.Item(248, 7).Value  = "Apio"
$ws.Cells.Item(248, 8).Value  = "Americana (o)"
$ws.Cells.Item(248, 9).Value  = "Primera"
$ws.Cells.Item(248, 10).Value = 50
$ws.Cells.Item(248, 11).Value = 9000
$ws.Cells.Item(248, 12).Value = 9000
$ws.Cells.Item(248, 13).Value = 9000
$ws.Cells.Item(248, 14).Value = "`$/docena de matas"
$ws.Cells.Item(248, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(248, 16).Value = 1500
$ws.Cells.Item(248, 17).Value = 6
$ws.Cells.Item(248, 18).Value = "Hortaliza"
